# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" '22.455.83'
Set-TextValue $ws "E2" '  +0.14%  '
Set-TextValue $ws "D3" '1.573.12'
Set-TextValue $ws "E4" '  -0.07%  '
Set-TextValue $ws "E5" '  -0.16%  '
Set-TextValue $ws "D6" '291.13'
Set-TextValue $ws "D7" '0.3739'
Set-TextValue $ws "E7" '  -0.79%  '
Set-TextValue $ws "E8" '  -0.04%  '
Set-TextValue $ws "D9" '0.3403'
Set-TextValue $ws "E9" '  -0.54%  '
Set-TextValue $ws "D10" '0.07557'
Set-TextValue $ws "E10" '  -1.50%  '
Set-TextValue $ws "D11" '1.140'
Set-TextValue $ws "E12" '  -0.03%  '
Set-TextValue $ws "D13" '21.33'
Set-TextValue $ws "E13" '  +0.18%  '
Set-TextValue $ws "E14" '  -0.26%  '
Set-TextValue $ws "D15" '6.940'
Set-TextValue $ws "E15" '  +0.09%  '
Set-TextValue $ws "D16" '1.569.98'
Set-TextValue $ws "E16" '  -0.69%  '
Set-TextValue $ws "D17" '0.00001121'
Set-TextValue $ws "E17" '  -1.44%  '
Set-TextValue $ws "E18" '  +0.60%  '
Set-TextValue $ws "D19" '0.06737'
Set-TextValue $ws "E19" '  -0.50%  '
Set-TextValue $ws "E20" '  -0.04%  '
Set-TextValue $ws "E21" '  +0.11%  '
Set-TextValue $ws "D22" '16.41'
Set-TextValue $ws "E22" '  -2.47%  '
Set-TextValue $ws "D23" '12.16'
Set-TextValue $ws "E23" '  +0.84%  '
Set-TextValue $ws "D24" '22.455.10'
Set-TextValue $ws "E24" '  +0.13%  '
Set-TextValue $ws "D25" '2.354'
Set-TextValue $ws "E25" '  -2.84%  '
Set-TextValue $ws "D26" '2.585'
Set-TextValue $ws "E26" '  -6.02%  '
Set-TextValue $ws "D27" '20.15'
Set-TextValue $ws "E27" '  -1.04%  '
Set-TextValue $ws "D28" '148.65'
Set-TextValue $ws "E28" '  +1.90%  '
Set-TextValue $ws "D29" '5.018'
Set-TextValue $ws "E29" '  -0.38%  '
Set-TextValue $ws "D30" '125.77'
Set-TextValue $ws "E30" '  -0.30%  '
Set-TextValue $ws "D31" '1.745.26'
Set-TextValue $ws "E31" '  -0.51%  '
Set-TextValue $ws "E32" '  +3.38%  '
Set-TextValue $ws "D33" '6.128'
Set-TextValue $ws "E33" '  -1.79%  '
Set-TextValue $ws "E34" '  -1.90%  '
Set-TextValue $ws "D35" '9.834'
Set-TextValue $ws "E35" '  -2.28%  '
Set-TextValue $ws "D36" '0.08394'
Set-TextValue $ws "E36" '  -2.11%  '
Set-TextValue $ws "D37" '1.376'
Set-TextValue $ws "E37" '  +2.54%  '
Set-TextValue $ws "E38" '  -3.90%  '
Set-TextValue $ws "E39" '  -0.86%  '
Set-TextValue $ws "D40" '0.06534'
Set-TextValue $ws "E40" '  -0.91%  '
Set-TextValue $ws "D41" '5.473'
Set-TextValue $ws "E41" '  +0.05%  '
Set-TextValue $ws "D42" '11.31'
Set-TextValue $ws "E42" '  -2.52%  '
Set-TextValue $ws "D43" '0.6271'
Set-TextValue $ws "E43" '  -3.05%  '
Set-TextValue $ws "E44" '  -0.08%  '
Set-TextValue $ws "D45" '13.97'
Set-TextValue $ws "E45" '  -1.00%  '
Set-TextValue $ws "D46" '3.808'
Set-TextValue $ws "E46" '  +0.20%  '
Set-TextValue $ws "D47" '0.5839'
Set-TextValue $ws "E47" '  -3.22%  '
Set-TextValue $ws "D48" '2.086'
Set-TextValue $ws "E48" '  -0.48%  '
Set-TextValue $ws "D49" '129.31'
Set-TextValue $ws "E49" '  +2.99%  '
Set-TextValue $ws "D50" '1.223'
Set-TextValue $ws "E50" '  -6.18%  '
Set-TextValue $ws "D51" '0.07328'
Set-TextValue $ws "E51" '  -0.05%  '